{"js": "// Replace the arithmetic problem text in each table cell, in document order,\n// while preserving existing run/paragraph formatting (only the digits/operator change).\nconst table = context.document.body.tables.getFirst();\n\nconst replacements = [\n  { row: 0, col: 0, oldText: \"94\u00f73=\", newText: \"33\u00f79=\" },\n  { row: 0, col: 1, oldText: \"48\u00f79=\", newText: \"63\u00f77=\" },\n  { row: 0, col: 2, oldText: \"62\u00f79=\", newText: \"55\u00f76=\" },\n  { row: 0, col: 3, oldText: \"68\u00f77=\", newText: \"80\u00f76=\" },\n  { row: 0, col: 4, oldText: \"18\u00f79=\", newText: \"98\u00f75=\" },\n  { row: 4, col: 0, oldText: \"21\u00f75=\", newText: \"19\u00f75=\" },\n  { row: 4, col: 1, oldText: \"84\u00f77=\", newText: \"35\u00f74=\" },\n  { row: 4, col: 2, oldText: \"68\u00f73=\", newText: \"66\u00f79=\" },\n  { row: 4, col: 3, oldText: \"66\u00f73=\", newText: \"85\u00f79=\" },\n  { row: 4, col: 4, oldText: \"36\u00f74=\", newText: \"66\u00f77=\" },\n  { row: 8, col: 0, oldText: \"34\u00f75=\", newText: \"40\u00f76=\" },\n  { row: 8, col: 1, oldText: \"65\u00f74=\", newText: \"52\u00f75=\" },\n  { row: 8, col: 2, oldText: \"21\u00f74=\", newText: \"14\u00f78=\" },\n  { row: 8, col: 3, oldText: \"12\u00f78=\", newText: \"87\u00f76=\" },\n  { row: 8, col: 4, oldText: \"48\u00f73=\", newText: \"12\u00f73=\" },\n  { row: 12, col: 0, oldText: \"70\u00f74=\", newText: \"68\u00f77=\" },\n  { row: 12, col: 1, oldText: \"89\u00f79=\", newText: \"52\u00f74=\" },\n  { row: 12, col: 2, oldText: \"13\u00f76=\", newText: \"23\u00f76=\" },\n  { row: 12, col: 3, oldText: \"33\u00f76=\", newText: \"20\u00f74=\" },\n  { row: 12, col: 4, oldText: \"54\u00f76=\", newText: \"24\u00f75=\" },\n  { row: 16, col: 0, oldText: \"50\u00f76=\", newText: \"63\u00f78=\" },\n  { row: 16, col: 1, oldText: \"33\u00f72=\", newText: \"49\u00f78=\" },\n  { row: 16, col: 2, oldText: \"70\u00f76=\", newText: \"56\u00f73=\" },\n  { row: 16, col: 3, oldText: \"84\u00f75=\", newText: \"53\u00f74=\" },\n  { row: 16, col: 4, oldText: \"25\u00f73=\", newText: \"54\u00f73=\" },\n];\n\nfor (const r of replacements) {\n  const cell = table.getCell(r.row, r.col);\n  const found = cell.body.search(r.oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n  if (found.items.length > 0) {\n    found.items[0].insertText(r.newText, Word.InsertLocation.replace);\n  } else {\n    // Fallback: replace the whole cell text if the exact match wasn't found.\n    cell.body.insertText(r.newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the arithmetic problem text in each table cell, in document order,\n# while preserving existing run/paragraph formatting (only the digits/operator change).\n$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n\n$replacements = @(\n    @{ Row = 1; Col = 1; OldText = \"94\u00f73=\"; NewText = \"33\u00f79=\" }\n    @{ Row = 1; Col = 2; OldText = \"48\u00f79=\"; NewText = \"63\u00f77=\" }\n    @{ Row = 1; Col = 3; OldText = \"62\u00f79=\"; NewText = \"55\u00f76=\" }\n    @{ Row = 1; Col = 4; OldText = \"68\u00f77=\"; NewText = \"80\u00f76=\" }\n    @{ Row = 1; Col = 5; OldText = \"18\u00f79=\"; NewText = \"98\u00f75=\" }\n    @{ Row = 5; Col = 1; OldText = \"21\u00f75=\"; NewText = \"19\u00f75=\" }\n    @{ Row = 5; Col = 2; OldText = \"84\u00f77=\"; NewText = \"35\u00f74=\" }\n    @{ Row = 5; Col = 3; OldText = \"68\u00f73=\"; NewText = \"66\u00f79=\" }\n    @{ Row = 5; Col = 4; OldText = \"66\u00f73=\"; NewText = \"85\u00f79=\" }\n    @{ Row = 5; Col = 5; OldText = \"36\u00f74=\"; NewText = \"66\u00f77=\" }\n    @{ Row = 9; Col = 1; OldText = \"34\u00f75=\"; NewText = \"40\u00f76=\" }\n    @{ Row = 9; Col = 2; OldText = \"65\u00f74=\"; NewText = \"52\u00f75=\" }\n    @{ Row = 9; Col = 3; OldText = \"21\u00f74=\"; NewText = \"14\u00f78=\" }\n    @{ Row = 9; Col = 4; OldText = \"12\u00f78=\"; NewText = \"87\u00f76=\" }\n    @{ Row = 9; Col = 5; OldText = \"48\u00f73=\"; NewText = \"12\u00f73=\" }\n    @{ Row = 13; Col = 1; OldText = \"70\u00f74=\"; NewText = \"68\u00f77=\" }\n    @{ Row = 13; Col = 2; OldText = \"89\u00f79=\"; NewText = \"52\u00f74=\" }\n    @{ Row = 13; Col = 3; OldText = \"13\u00f76=\"; NewText = \"23\u00f76=\" }\n    @{ Row = 13; Col = 4; OldText = \"33\u00f76=\"; NewText = \"20\u00f74=\" }\n    @{ Row = 13; Col = 5; OldText = \"54\u00f76=\"; NewText = \"24\u00f75=\" }\n    @{ Row = 17; Col = 1; OldText = \"50\u00f76=\"; NewText = \"63\u00f78=\" }\n    @{ Row = 17; Col = 2; OldText = \"33\u00f72=\"; NewText = \"49\u00f78=\" }\n    @{ Row = 17; Col = 3; OldText = \"70\u00f76=\"; NewText = \"56\u00f73=\" }\n    @{ Row = 17; Col = 4; OldText = \"84\u00f75=\"; NewText = \"53\u00f74=\" }\n    @{ Row = 17; Col = 5; OldText = \"25\u00f73=\"; NewText = \"54\u00f73=\" }\n)\n\nforeach ($r in $replacements) {\n    $cell = $table.Cell($r.Row, $r.Col)\n    $rng = $cell.Range\n    # wdFindContinue=1, wdReplaceOne=1 (replace just the first/only match in this cell)\n    $found = $rng.Find.Execute($r.OldText, $false, $false, $false, $false, $false, $true, 1, $false, $r.NewText, 1)\n    if (-not $found) {\n        Write-Output (\"WARN: not found -> row=\" + $r.Row + \" col=\" + $r.Col)\n    }\n}\n"}
